$d = $word.ActiveDocument

# --- 1) Remove the old "_GoBack" bookmark first (it currently sits right
#        after the "...order mentioned." run, further down the document).
#        Doing this before inserting the new bookmark avoids any ambiguity
#        about which same-named bookmark a later lookup would resolve to. ---
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# --- 2) Split the opening line's single run into two runs with new wording,
#        and add the "_GoBack" bookmark here (right after the new text). ---
$firstPara = $d.Paragraphs(1).Range
$oldText = "Message sequence for operations in Btalk:"
$firstLen = $oldText.Length
$target = $d.Range($firstPara.Start, $firstPara.Start + $firstLen)

if ($target.Text -ne $oldText) {
    throw "Unexpected content in first paragraph: [$($target.Text)]"
}

$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Message s</w:t></w:r><w:r><w:t>equence for operations in Btalk (unfinished version)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($newXml)

Write-Output "done"
